$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2 value from 1 to 2 (last control point count)
$ws.Range("B2").Value = 2

# Remove the last row (row 4: A4=2, B4=1) - "last cases" trimmed from cluster control points
$ws.Rows("4:4").Delete()
